# Add a header row (name / item / rate per unit) to Sheet1, size the
# "rate per unit" column to fit, and leave the selection on D1 - matching
# what Excel records after typing the three header cells and widening
# column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name "
$ws.Range("B1").Value = "item"
$ws.Range("C1").Value = "rate per unit"

# Widen column C to fit "rate per unit" (stored column width ends up at 12
# characters; ColumnWidth and the on-disk <col width> differ by the fixed
# ~0.83 char padding Excel applies, so back that out here).
$ws.Columns.Item(3).ColumnWidth = 67/6

$ws.Range("D1").Select()
